# "Update UI: Login + Sign up"
# Mark the three UI-design work items (rows 5-7, column E "Muc Do Hoan Thanh")
# with their actual completion percentage, matching the existing style
# already used by the sibling cells E2:E4 (percentage number format,
# centered both ways), then leave the selection where the user's cursor
# ended up (E7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 - "Thiet ke man hinh: Tao de xuat / Tao bai viet / ..." -> done
# Row 6 - "Thiet ke man hinh: Dang nhap / Dang ky / ..."         -> done
# Row 7 - "Thiet ke man hinh: Dang nhap / Dang ky tai khoan /..." -> not done yet

$progressRange = $ws.Range("E5:E7")
$progressRange.NumberFormat = "0%"
$progressRange.HorizontalAlignment = -4108  # xlCenter
$progressRange.VerticalAlignment = -4108    # xlCenter

$ws.Range("E5").Value = 1
$ws.Range("E6").Value = 1
$ws.Range("E7").Value = 0

# Leave the active selection on E7, matching where editing finished.
[void]$ws.Range("E7").Select()
